# Added Hungary FC Test data
# For each market sheet (UK, Spain, Hungary, Italy) insert the two new
# accessory codes "MX-BBX" and "MX-DPBX" into column A, just above the
# existing "Wg"/"Accessories" sentinel rows at the bottom of the list.

$wb = $excel.ActiveWorkbook

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Push the last two rows (the "Wg" / "Accessories" sentinel rows) down
    # by inserting two fresh rows above them, at row 20.
    $ws.Rows.Item(20).Insert()
    $ws.Rows.Item(20).Insert()

    # Copy the formatting (border/style) of the row above (row 19, which
    # carries the list-item style) onto the two new rows so they match the
    # rest of the list instead of being left unstyled.
    $ws.Range("A19").Copy()
    $ws.Range("A20:A21").PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    $ws.Range("A20").Value = "MX-BBX"
    $ws.Range("A21").Value = "MX-DPBX"
}

# Update the sheet-level selections to what they were left at after the
# edit, and make "Hungary" (the 3rd sheet) the active tab.
$wsUK = $wb.Worksheets.Item(1)
$wsUK.Activate()
$wsUK.Range("A7:A23").Select()

$wsSpain = $wb.Worksheets.Item(2)
$wsSpain.Activate()
$wsSpain.Range("A7:A23").Select()

$wsItaly = $wb.Worksheets.Item(4)
$wsItaly.Activate()
$wsItaly.Range("A13").Select()

$wsHungary = $wb.Worksheets.Item(3)
$wsHungary.Activate()
$wsHungary.Range("A10").Select()
